$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Status column: mark the "Insert" steps (rows 7,9,11,13,15) as done ("!")
# instead of pending ("?").
foreach ($r in 7,9,11,13,15) {
    $ws.Cells.Item($r, 6).Value = "!"
}

# Fix the typo in the table name for the "venda_Itens" task.
$ws.Range("B16").Value = "Tabela venda_Itens"

# Leave the selection on the cell that was last edited.
[void]$ws.Range("B16").Select()
